$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing values in column B
$ws.Range("B1").Value = 228
$ws.Range("B3").Value = 0
$ws.Range("B4").Value = 346
$ws.Range("B7").Value = 86
$ws.Range("B9").Value = 85

# Add new rows 13 and 14 (column A holds large Discord snowflake IDs stored as text)
$ws.Range("A13").NumberFormat = "@"
$ws.Range("A13").Value = "609604172349964328"
$ws.Range("B13").Value = 318
$ws.Range("C13").Value = "귤님#2613"

$ws.Range("A14").NumberFormat = "@"
$ws.Range("A14").Value = "462469630347182080"
$ws.Range("B14").Value = 68
$ws.Range("C14").Value = "☔ bow 6#0739"
